$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2 through 6:
#  - Column B ("Did Harvest Occur?"): "Yes" -> "No"
#  - Column F ("Species"): "Na" -> "" (blank)
#  - Column J ("Unknown Sex Count"): 1 -> 0
for ($row = 2; $row -le 6; $row++) {
    $ws.Range("B$row").Value = "No"
    $ws.Range("F$row").Value = ""
    $ws.Range("J$row").Value = 0
}
